$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 45931
$ws.Cells.Item(2, 2).Value = 0.02
$ws.Cells.Item(2, 3).Value = 2.099
$ws.Cells.Item(2, 4).Value = 0
$ws.Cells.Item(2, 5).Value = 0
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = "01.10.20251"

$ws.Cells.Item(3, 1).Value = 45931.01041666666
$ws.Cells.Item(3, 2).Value = 0.08699999999999999
$ws.Cells.Item(3, 3).Value = 1.343
$ws.Cells.Item(3, 4).Value = 0
$ws.Cells.Item(3, 5).Value = 0
$ws.Cells.Item(3, 6).Value = 2
$ws.Cells.Item(3, 7).Value = "01.10.20252"

$ws.Cells.Item(4, 1).Value = 45931.02083333334
$ws.Cells.Item(4, 2).Value = 0.082
$ws.Cells.Item(4, 3).Value = 2.967
$ws.Cells.Item(4, 4).Value = 0
$ws.Cells.Item(4, 5).Value = 0
$ws.Cells.Item(4, 6).Value = 3
$ws.Cells.Item(4, 7).Value = "01.10.20253"

$ws.Cells.Item(5, 1).Value = 45931.03125
$ws.Cells.Item(5, 2).Value = 0.149
$ws.Cells.Item(5, 3).Value = 6.195
$ws.Cells.Item(5, 4).Value = 0
$ws.Cells.Item(5, 5).Value = 0
$ws.Cells.Item(5, 6).Value = 4
$ws.Cells.Item(5, 7).Value = "01.10.20254"

$ws.Cells.Item(6, 1).Value = 45931.04166666666
$ws.Cells.Item(6, 2).Value = 0
$ws.Cells.Item(6, 3).Value = 1.388
$ws.Cells.Item(6, 4).Value = 0
$ws.Cells.Item(6, 5).Value = 0
$ws.Cells.Item(6, 6).Value = 5
$ws.Cells.Item(6, 7).Value = "01.10.20255"

$ws.Cells.Item(7, 1).Value = 45931.05208333334
$ws.Cells.Item(7, 2).Value = 0.017
$ws.Cells.Item(7, 3).Value = 2.756
$ws.Cells.Item(7, 4).Value = 0
$ws.Cells.Item(7, 5).Value = 0
$ws.Cells.Item(7, 6).Value = 6
$ws.Cells.Item(7, 7).Value = "01.10.20256"

$ws.Cells.Item(8, 1).Value = 45931.0625
$ws.Cells.Item(8, 2).Value = 0
$ws.Cells.Item(8, 3).Value = 4.665
$ws.Cells.Item(8, 4).Value = 0
$ws.Cells.Item(8, 5).Value = 0
$ws.Cells.Item(8, 6).Value = 7
$ws.Cells.Item(8, 7).Value = "01.10.20257"

$ws.Cells.Item(9, 1).Value = 45931.07291666666
$ws.Cells.Item(9, 2).Value = 0
$ws.Cells.Item(9, 3).Value = 12.601
$ws.Cells.Item(9, 4).Value = 0
$ws.Cells.Item(9, 5).Value = 0
$ws.Cells.Item(9, 6).Value = 8
$ws.Cells.Item(9, 7).Value = "01.10.20258"

$ws.Cells.Item(10, 1).Value = 45931.08333333334
$ws.Cells.Item(10, 2).Value = 0
$ws.Cells.Item(10, 3).Value = 18.553
$ws.Cells.Item(10, 4).Value = 0
$ws.Cells.Item(10, 5).Value = 0
$ws.Cells.Item(10, 6).Value = 9
$ws.Cells.Item(10, 7).Value = "01.10.20259"

$ws.Cells.Item(11, 1).Value = 45931.09375
$ws.Cells.Item(11, 2).Value = 0
$ws.Cells.Item(11, 3).Value = 13.939
$ws.Cells.Item(11, 4).Value = 0
$ws.Cells.Item(11, 5).Value = 0
$ws.Cells.Item(11, 6).Value = 10
$ws.Cells.Item(11, 7).Value = "01.10.202510"

$ws.Cells.Item(12, 1).Value = 45931.10416666666
$ws.Cells.Item(12, 2).Value = 0
$ws.Cells.Item(12, 3).Value = 21.859
$ws.Cells.Item(12, 4).Value = 0
$ws.Cells.Item(12, 5).Value = 0
$ws.Cells.Item(12, 6).Value = 11
$ws.Cells.Item(12, 7).Value = "01.10.202511"

$ws.Cells.Item(13, 1).Value = 45931.11458333334
$ws.Cells.Item(13, 2).Value = 0
$ws.Cells.Item(13, 3).Value = 15.319
$ws.Cells.Item(13, 4).Value = 0
$ws.Cells.Item(13, 5).Value = 0
$ws.Cells.Item(13, 6).Value = 12
$ws.Cells.Item(13, 7).Value = "01.10.202512"

$ws.Cells.Item(14, 1).Value = 45931.125
$ws.Cells.Item(14, 2).Value = 0
$ws.Cells.Item(14, 3).Value = 8.343999999999999
$ws.Cells.Item(14, 4).Value = 0
$ws.Cells.Item(14, 5).Value = 0
$ws.Cells.Item(14, 6).Value = 13
$ws.Cells.Item(14, 7).Value = "01.10.202513"

$ws.Cells.Item(15, 1).Value = 45931.13541666666
$ws.Cells.Item(15, 2).Value = 0
$ws.Cells.Item(15, 3).Value = 8.818
$ws.Cells.Item(15, 4).Value = 0
$ws.Cells.Item(15, 5).Value = 0
$ws.Cells.Item(15, 6).Value = 14
$ws.Cells.Item(15, 7).Value = "01.10.202514"

$ws.Cells.Item(16, 1).Value = 45931.14583333334
$ws.Cells.Item(16, 2).Value = 0.013
$ws.Cells.Item(16, 3).Value = 4.677
$ws.Cells.Item(16, 4).Value = 0
$ws.Cells.Item(16, 5).Value = 0
$ws.Cells.Item(16, 6).Value = 15
$ws.Cells.Item(16, 7).Value = "01.10.202515"

$ws.Cells.Item(17, 1).Value = 45931.15625
$ws.Cells.Item(17, 2).Value = 0
$ws.Cells.Item(17, 3).Value = 3.812
$ws.Cells.Item(17, 4).Value = 0
$ws.Cells.Item(17, 5).Value = 0
$ws.Cells.Item(17, 6).Value = 16
$ws.Cells.Item(17, 7).Value = "01.10.202516"

$ws.Cells.Item(18, 1).Value = 45931.16666666666
$ws.Cells.Item(18, 2).Value = 4.493
$ws.Cells.Item(18, 3).Value = 0.362
$ws.Cells.Item(18, 4).Value = 0
$ws.Cells.Item(18, 5).Value = 0
$ws.Cells.Item(18, 6).Value = 17
$ws.Cells.Item(18, 7).Value = "01.10.202517"

$ws.Cells.Item(19, 1).Value = 45931.17708333334
$ws.Cells.Item(19, 2).Value = 4.433
$ws.Cells.Item(19, 3).Value = 0
$ws.Cells.Item(19, 4).Value = 0
$ws.Cells.Item(19, 5).Value = 0
$ws.Cells.Item(19, 6).Value = 18
$ws.Cells.Item(19, 7).Value = "01.10.202518"

$ws.Cells.Item(20, 1).Value = 45931.1875
$ws.Cells.Item(20, 2).Value = 9.034000000000001
$ws.Cells.Item(20, 3).Value = 0
$ws.Cells.Item(20, 4).Value = 0
$ws.Cells.Item(20, 5).Value = 0
$ws.Cells.Item(20, 6).Value = 19
$ws.Cells.Item(20, 7).Value = "01.10.202519"

$ws.Cells.Item(21, 1).Value = 45931.19791666666
$ws.Cells.Item(21, 2).Value = 1.039
$ws.Cells.Item(21, 3).Value = 0.597
$ws.Cells.Item(21, 4).Value = 0
$ws.Cells.Item(21, 5).Value = 0
$ws.Cells.Item(21, 6).Value = 20
$ws.Cells.Item(21, 7).Value = "01.10.202520"

$ws.Cells.Item(22, 1).Value = 45931.20833333334
$ws.Cells.Item(22, 2).Value = 0.021
$ws.Cells.Item(22, 3).Value = 1.178
$ws.Cells.Item(22, 4).Value = 0
$ws.Cells.Item(22, 5).Value = 0
$ws.Cells.Item(22, 6).Value = 21
$ws.Cells.Item(22, 7).Value = "01.10.202521"

$ws.Cells.Item(23, 1).Value = 45931.21875
$ws.Cells.Item(23, 2).Value = 0
$ws.Cells.Item(23, 3).Value = 10.938
$ws.Cells.Item(23, 4).Value = 0
$ws.Cells.Item(23, 5).Value = 0
$ws.Cells.Item(23, 6).Value = 22
$ws.Cells.Item(23, 7).Value = "01.10.202522"

$ws.Cells.Item(24, 1).Value = 45931.22916666666
$ws.Cells.Item(24, 2).Value = 0.37
$ws.Cells.Item(24, 3).Value = 2.106
$ws.Cells.Item(24, 4).Value = 0
$ws.Cells.Item(24, 5).Value = 0
$ws.Cells.Item(24, 6).Value = 23
$ws.Cells.Item(24, 7).Value = "01.10.202523"

$ws.Cells.Item(25, 1).Value = 45931.23958333334
$ws.Cells.Item(25, 2).Value = 0.585
$ws.Cells.Item(25, 3).Value = 4.713
$ws.Cells.Item(25, 4).Value = 0
$ws.Cells.Item(25, 5).Value = 0
$ws.Cells.Item(25, 6).Value = 24
$ws.Cells.Item(25, 7).Value = "01.10.202524"

$ws.Cells.Item(26, 1).Value = 45931.25
$ws.Cells.Item(26, 2).Value = 1.101
$ws.Cells.Item(26, 3).Value = 1.79
$ws.Cells.Item(26, 4).Value = 0
$ws.Cells.Item(26, 5).Value = 0
$ws.Cells.Item(26, 6).Value = 25
$ws.Cells.Item(26, 7).Value = "01.10.202525"

$ws.Cells.Item(27, 1).Value = 45931.26041666666
$ws.Cells.Item(27, 2).Value = 0.722
$ws.Cells.Item(27, 3).Value = 0.14
$ws.Cells.Item(27, 4).Value = 0
$ws.Cells.Item(27, 5).Value = 0
$ws.Cells.Item(27, 6).Value = 26
$ws.Cells.Item(27, 7).Value = "01.10.202526"

$ws.Cells.Item(28, 1).Value = 45931.27083333334
$ws.Cells.Item(28, 2).Value = 0.607
$ws.Cells.Item(28, 3).Value = 0.045
$ws.Cells.Item(28, 4).Value = 0
$ws.Cells.Item(28, 5).Value = 0
$ws.Cells.Item(28, 6).Value = 27
$ws.Cells.Item(28, 7).Value = "01.10.202527"

$ws.Cells.Item(29, 1).Value = 45931.28125
$ws.Cells.Item(29, 2).Value = 0.283
$ws.Cells.Item(29, 3).Value = 0.054
$ws.Cells.Item(29, 4).Value = 0
$ws.Cells.Item(29, 5).Value = 0
$ws.Cells.Item(29, 6).Value = 28
$ws.Cells.Item(29, 7).Value = "01.10.202528"

$ws.Cells.Item(30, 1).Value = 45931.29166666666
$ws.Cells.Item(30, 2).Value = 5.397
$ws.Cells.Item(30, 3).Value = 0.026
$ws.Cells.Item(30, 4).Value = 0
$ws.Cells.Item(30, 5).Value = 0
$ws.Cells.Item(30, 6).Value = 29
$ws.Cells.Item(30, 7).Value = "01.10.202529"

$ws.Cells.Item(31, 1).Value = 45931.30208333334
$ws.Cells.Item(31, 2).Value = 1.248
$ws.Cells.Item(31, 3).Value = 0.016
$ws.Cells.Item(31, 4).Value = 0
$ws.Cells.Item(31, 5).Value = 0
$ws.Cells.Item(31, 6).Value = 30
$ws.Cells.Item(31, 7).Value = "01.10.202530"

$ws.Cells.Item(32, 1).Value = 45931.3125
$ws.Cells.Item(32, 2).Value = 1.411
$ws.Cells.Item(32, 3).Value = 1.302
$ws.Cells.Item(32, 4).Value = 0
$ws.Cells.Item(32, 5).Value = 0
$ws.Cells.Item(32, 6).Value = 31
$ws.Cells.Item(32, 7).Value = "01.10.202531"

$ws.Cells.Item(33, 1).Value = 45931.32291666666
$ws.Cells.Item(33, 2).Value = 0
$ws.Cells.Item(33, 3).Value = 3.204
$ws.Cells.Item(33, 4).Value = 0
$ws.Cells.Item(33, 5).Value = 0
$ws.Cells.Item(33, 6).Value = 32
$ws.Cells.Item(33, 7).Value = "01.10.202532"

$ws.Cells.Item(34, 1).Value = 45931.33333333334
$ws.Cells.Item(34, 2).Value = 3.726
$ws.Cells.Item(34, 3).Value = 0.401
$ws.Cells.Item(34, 4).Value = 0
$ws.Cells.Item(34, 5).Value = 0
$ws.Cells.Item(34, 6).Value = 33
$ws.Cells.Item(34, 7).Value = "01.10.202533"

$ws.Cells.Item(35, 1).Value = 45931.34375
$ws.Cells.Item(35, 2).Value = 0.09
$ws.Cells.Item(35, 3).Value = 0.125
$ws.Cells.Item(35, 4).Value = 0
$ws.Cells.Item(35, 5).Value = 0
$ws.Cells.Item(35, 6).Value = 34
$ws.Cells.Item(35, 7).Value = "01.10.202534"

$ws.Cells.Item(36, 1).Value = 45931.35416666666
$ws.Cells.Item(36, 2).Value = 0.091
$ws.Cells.Item(36, 3).Value = 7.242
$ws.Cells.Item(36, 4).Value = 0
$ws.Cells.Item(36, 5).Value = 0
$ws.Cells.Item(36, 6).Value = 35
$ws.Cells.Item(36, 7).Value = "01.10.202535"

$ws.Cells.Item(37, 1).Value = 45931.36458333334
$ws.Cells.Item(37, 2).Value = 0
$ws.Cells.Item(37, 3).Value = 20.918
$ws.Cells.Item(37, 4).Value = 0
$ws.Cells.Item(37, 5).Value = 0
$ws.Cells.Item(37, 6).Value = 36
$ws.Cells.Item(37, 7).Value = "01.10.202536"

$ws.Cells.Item(38, 1).Value = 45931.375
$ws.Cells.Item(38, 2).Value = 0.004
$ws.Cells.Item(38, 3).Value = 10.637
$ws.Cells.Item(38, 4).Value = 0
$ws.Cells.Item(38, 5).Value = 0
$ws.Cells.Item(38, 6).Value = 37
$ws.Cells.Item(38, 7).Value = "01.10.202537"

$ws.Cells.Item(39, 1).Value = 45931.38541666666
$ws.Cells.Item(39, 2).Value = 0
$ws.Cells.Item(39, 3).Value = 17.613
$ws.Cells.Item(39, 4).Value = 0
$ws.Cells.Item(39, 5).Value = 0
$ws.Cells.Item(39, 6).Value = 38
$ws.Cells.Item(39, 7).Value = "01.10.202538"

$ws.Cells.Item(40, 1).Value = 45931.39583333334
$ws.Cells.Item(40, 2).Value = 0
$ws.Cells.Item(40, 3).Value = 9.992000000000001
$ws.Cells.Item(40, 4).Value = 0
$ws.Cells.Item(40, 5).Value = 0
$ws.Cells.Item(40, 6).Value = 39
$ws.Cells.Item(40, 7).Value = "01.10.202539"

$ws.Cells.Item(41, 1).Value = 45931.40625
$ws.Cells.Item(41, 2).Value = 0
$ws.Cells.Item(41, 3).Value = 21.155
$ws.Cells.Item(41, 4).Value = 0
$ws.Cells.Item(41, 5).Value = 0
$ws.Cells.Item(41, 6).Value = 40
$ws.Cells.Item(41, 7).Value = "01.10.202540"

$ws.Cells.Item(42, 1).Value = 45931.41666666666
$ws.Cells.Item(42, 2).Value = 0
$ws.Cells.Item(42, 3).Value = 26.094
$ws.Cells.Item(42, 4).Value = 0
$ws.Cells.Item(42, 5).Value = 0
$ws.Cells.Item(42, 6).Value = 41
$ws.Cells.Item(42, 7).Value = "01.10.202541"

$ws.Cells.Item(43, 1).Value = 45931.42708333334
$ws.Cells.Item(43, 2).Value = 0
$ws.Cells.Item(43, 3).Value = 39.129
$ws.Cells.Item(43, 4).Value = 0
$ws.Cells.Item(43, 5).Value = 0
$ws.Cells.Item(43, 6).Value = 42
$ws.Cells.Item(43, 7).Value = "01.10.202542"

$ws.Cells.Item(44, 1).Value = 45931.4375
$ws.Cells.Item(44, 2).Value = 0
$ws.Cells.Item(44, 3).Value = 36.674
$ws.Cells.Item(44, 4).Value = 0
$ws.Cells.Item(44, 5).Value = 0
$ws.Cells.Item(44, 6).Value = 43
$ws.Cells.Item(44, 7).Value = "01.10.202543"

$ws.Cells.Item(45, 1).Value = 45931.44791666666
$ws.Cells.Item(45, 2).Value = 0
$ws.Cells.Item(45, 3).Value = 27.657
$ws.Cells.Item(45, 4).Value = 0
$ws.Cells.Item(45, 5).Value = 25
$ws.Cells.Item(45, 6).Value = 44
$ws.Cells.Item(45, 7).Value = "01.10.202544"

$ws.Cells.Item(46, 1).Value = 45931.45833333334
$ws.Cells.Item(46, 2).Value = 0
$ws.Cells.Item(46, 3).Value = 38.035
$ws.Cells.Item(46, 4).Value = 0
$ws.Cells.Item(46, 5).Value = 25
$ws.Cells.Item(46, 6).Value = 45
$ws.Cells.Item(46, 7).Value = "01.10.202545"

$ws.Cells.Item(47, 1).Value = 45931.46875
$ws.Cells.Item(47, 2).Value = 0
$ws.Cells.Item(47, 3).Value = 23.227
$ws.Cells.Item(47, 4).Value = 0
$ws.Cells.Item(47, 5).Value = 25
$ws.Cells.Item(47, 6).Value = 46
$ws.Cells.Item(47, 7).Value = "01.10.202546"

$ws.Cells.Item(48, 1).Value = 45931.47916666666
$ws.Cells.Item(48, 2).Value = 0
$ws.Cells.Item(48, 3).Value = 8.281000000000001
$ws.Cells.Item(48, 4).Value = 0
$ws.Cells.Item(48, 5).Value = 44.75
$ws.Cells.Item(48, 6).Value = 47
$ws.Cells.Item(48, 7).Value = "01.10.202547"

$ws.Cells.Item(49, 1).Value = 45931.48958333334
$ws.Cells.Item(49, 2).Value = 0
$ws.Cells.Item(49, 3).Value = 9.685
$ws.Cells.Item(49, 4).Value = 0
$ws.Cells.Item(49, 5).Value = 44.75
$ws.Cells.Item(49, 6).Value = 48
$ws.Cells.Item(49, 7).Value = "01.10.202548"

$ws.Cells.Item(50, 1).Value = 45931.5
$ws.Cells.Item(50, 2).Value = 0
$ws.Cells.Item(50, 3).Value = 32.053
$ws.Cells.Item(50, 4).Value = 0
$ws.Cells.Item(50, 5).Value = 44.75
$ws.Cells.Item(50, 6).Value = 49
$ws.Cells.Item(50, 7).Value = "01.10.202549"

$ws.Cells.Item(51, 1).Value = 45931.51041666666
$ws.Cells.Item(51, 2).Value = 0
$ws.Cells.Item(51, 3).Value = 8.106999999999999
$ws.Cells.Item(51, 4).Value = 0
$ws.Cells.Item(51, 5).Value = 44.75
$ws.Cells.Item(51, 6).Value = 50
$ws.Cells.Item(51, 7).Value = "01.10.202550"

$ws.Cells.Item(52, 1).Value = 45931.52083333334
$ws.Cells.Item(52, 2).Value = 0
$ws.Cells.Item(52, 3).Value = 6.841
$ws.Cells.Item(52, 4).Value = 0
$ws.Cells.Item(52, 5).Value = 44.75
$ws.Cells.Item(52, 6).Value = 51
$ws.Cells.Item(52, 7).Value = "01.10.202551"

$ws.Cells.Item(53, 1).Value = 45931.53125
$ws.Cells.Item(53, 2).Value = 0.113
$ws.Cells.Item(53, 3).Value = 2.735
$ws.Cells.Item(53, 4).Value = 0
$ws.Cells.Item(53, 5).Value = 75
$ws.Cells.Item(53, 6).Value = 52
$ws.Cells.Item(53, 7).Value = "01.10.202552"

$ws.Cells.Item(54, 1).Value = 45931.54166666666
$ws.Cells.Item(54, 2).Value = 0.121
$ws.Cells.Item(54, 3).Value = 2.942
$ws.Cells.Item(54, 4).Value = 0
$ws.Cells.Item(54, 5).Value = 75
$ws.Cells.Item(54, 6).Value = 53
$ws.Cells.Item(54, 7).Value = "01.10.202553"

$ws.Cells.Item(55, 1).Value = 45931.55208333334
$ws.Cells.Item(55, 2).Value = 1.323
$ws.Cells.Item(55, 3).Value = 0.043
$ws.Cells.Item(55, 4).Value = 0
$ws.Cells.Item(55, 5).Value = 75
$ws.Cells.Item(55, 6).Value = 54
$ws.Cells.Item(55, 7).Value = "01.10.202554"

$ws.Cells.Item(56, 1).Value = 45931.5625
$ws.Cells.Item(56, 2).Value = 2.173
$ws.Cells.Item(56, 3).Value = 0.119
$ws.Cells.Item(56, 4).Value = 0
$ws.Cells.Item(56, 5).Value = 37.5
$ws.Cells.Item(56, 6).Value = 55
$ws.Cells.Item(56, 7).Value = "01.10.202555"

$ws.Cells.Item(57, 1).Value = 45931.57291666666
$ws.Cells.Item(57, 2).Value = 0.445
$ws.Cells.Item(57, 3).Value = 1.805
$ws.Cells.Item(57, 4).Value = 0
$ws.Cells.Item(57, 5).Value = 25
$ws.Cells.Item(57, 6).Value = 56
$ws.Cells.Item(57, 7).Value = "01.10.202556"

$ws.Cells.Item(58, 1).Value = 45931.58333333334
$ws.Cells.Item(58, 2).Value = 0
$ws.Cells.Item(58, 3).Value = 3.862
$ws.Cells.Item(58, 4).Value = 0
$ws.Cells.Item(58, 5).Value = 0
$ws.Cells.Item(58, 6).Value = 57
$ws.Cells.Item(58, 7).Value = "01.10.202557"

$ws.Cells.Item(59, 1).Value = 45931.59375
$ws.Cells.Item(59, 2).Value = 1.02
$ws.Cells.Item(59, 3).Value = 0.8139999999999999
$ws.Cells.Item(59, 4).Value = 0
$ws.Cells.Item(59, 5).Value = 0
$ws.Cells.Item(59, 6).Value = 58
$ws.Cells.Item(59, 7).Value = "01.10.202558"

$ws.Cells.Item(60, 1).Value = 45931.60416666666
$ws.Cells.Item(60, 2).Value = 0.134
$ws.Cells.Item(60, 3).Value = 1.433
$ws.Cells.Item(60, 4).Value = 0
$ws.Cells.Item(60, 5).Value = 0
$ws.Cells.Item(60, 6).Value = 59
$ws.Cells.Item(60, 7).Value = "01.10.202559"

$ws.Cells.Item(61, 1).Value = 45931.61458333334
$ws.Cells.Item(61, 2).Value = 0.09
$ws.Cells.Item(61, 3).Value = 0.503
$ws.Cells.Item(61, 4).Value = 0
$ws.Cells.Item(61, 5).Value = 0
$ws.Cells.Item(61, 6).Value = 60
$ws.Cells.Item(61, 7).Value = "01.10.202560"

$ws.Cells.Item(62, 1).Value = 45931.625
$ws.Cells.Item(62, 2).Value = 0.008
$ws.Cells.Item(62, 3).Value = 1.594
$ws.Cells.Item(62, 4).Value = 0
$ws.Cells.Item(62, 5).Value = 0
$ws.Cells.Item(62, 6).Value = 61
$ws.Cells.Item(62, 7).Value = "01.10.202561"

$ws.Cells.Item(63, 1).Value = 45931.63541666666
$ws.Cells.Item(63, 2).Value = 0.202
$ws.Cells.Item(63, 3).Value = 0.093
$ws.Cells.Item(63, 4).Value = 0
$ws.Cells.Item(63, 5).Value = 0
$ws.Cells.Item(63, 6).Value = 62
$ws.Cells.Item(63, 7).Value = "01.10.202562"

$ws.Cells.Item(64, 1).Value = 45931.64583333334
$ws.Cells.Item(64, 2).Value = 2.398
$ws.Cells.Item(64, 3).Value = 0.089
$ws.Cells.Item(64, 4).Value = 0
$ws.Cells.Item(64, 5).Value = 0
$ws.Cells.Item(64, 6).Value = 63
$ws.Cells.Item(64, 7).Value = "01.10.202563"

$ws.Cells.Item(65, 1).Value = 45931.65625
$ws.Cells.Item(65, 2).Value = 2.265
$ws.Cells.Item(65, 3).Value = 0
$ws.Cells.Item(65, 4).Value = 0
$ws.Cells.Item(65, 5).Value = 0
$ws.Cells.Item(65, 6).Value = 64
$ws.Cells.Item(65, 7).Value = "01.10.202564"

$ws.Cells.Item(66, 1).Value = 45931.66666666666
$ws.Cells.Item(66, 2).Value = 0.179
$ws.Cells.Item(66, 3).Value = 0.216
$ws.Cells.Item(66, 4).Value = 0
$ws.Cells.Item(66, 5).Value = 0
$ws.Cells.Item(66, 6).Value = 65
$ws.Cells.Item(66, 7).Value = "01.10.202565"

$ws.Cells.Item(67, 1).Value = 45931.67708333334
$ws.Cells.Item(67, 2).Value = 0.105
$ws.Cells.Item(67, 3).Value = 0.303
$ws.Cells.Item(67, 4).Value = 0
$ws.Cells.Item(67, 5).Value = 0
$ws.Cells.Item(67, 6).Value = 66
$ws.Cells.Item(67, 7).Value = "01.10.202566"

$ws.Cells.Item(68, 1).Value = 45931.6875
$ws.Cells.Item(68, 2).Value = 0.007
$ws.Cells.Item(68, 3).Value = 0.66
$ws.Cells.Item(68, 4).Value = 0
$ws.Cells.Item(68, 5).Value = 0
$ws.Cells.Item(68, 6).Value = 67
$ws.Cells.Item(68, 7).Value = "01.10.202567"

$ws.Cells.Item(69, 1).Value = 45931.69791666666
$ws.Cells.Item(69, 2).Value = 0.006
$ws.Cells.Item(69, 3).Value = 2.989
$ws.Cells.Item(69, 4).Value = 0
$ws.Cells.Item(69, 5).Value = 0
$ws.Cells.Item(69, 6).Value = 68
$ws.Cells.Item(69, 7).Value = "01.10.202568"

$ws.Cells.Item(70, 1).Value = 45931.70833333334
$ws.Cells.Item(70, 2).Value = 0
$ws.Cells.Item(70, 3).Value = 12.056
$ws.Cells.Item(70, 4).Value = 0
$ws.Cells.Item(70, 5).Value = 0
$ws.Cells.Item(70, 6).Value = 69
$ws.Cells.Item(70, 7).Value = "01.10.202569"

$ws.Cells.Item(71, 1).Value = 45931.71875
$ws.Cells.Item(71, 2).Value = 0.195
$ws.Cells.Item(71, 3).Value = 8.686
$ws.Cells.Item(71, 4).Value = 0
$ws.Cells.Item(71, 5).Value = 87.5
$ws.Cells.Item(71, 6).Value = 70
$ws.Cells.Item(71, 7).Value = "01.10.202570"

$ws.Cells.Item(72, 1).Value = 45931.72916666666
$ws.Cells.Item(72, 2).Value = 0.886
$ws.Cells.Item(72, 3).Value = 0.283
$ws.Cells.Item(72, 4).Value = 0
$ws.Cells.Item(72, 5).Value = 68.5
$ws.Cells.Item(72, 6).Value = 71
$ws.Cells.Item(72, 7).Value = "01.10.202571"

$ws.Cells.Item(73, 1).Value = 45931.73958333334
$ws.Cells.Item(73, 2).Value = 0.105
$ws.Cells.Item(73, 3).Value = 2.079
$ws.Cells.Item(73, 4).Value = 0
$ws.Cells.Item(73, 5).Value = 68.5
$ws.Cells.Item(73, 6).Value = 72
$ws.Cells.Item(73, 7).Value = "01.10.202572"

$ws.Cells.Item(74, 1).Value = 45931.75
$ws.Cells.Item(74, 2).Value = 0
$ws.Cells.Item(74, 3).Value = 8.77
$ws.Cells.Item(74, 4).Value = 0
$ws.Cells.Item(74, 5).Value = 75
$ws.Cells.Item(74, 6).Value = 73
$ws.Cells.Item(74, 7).Value = "01.10.202573"

$ws.Cells.Item(75, 1).Value = 45931.76041666666
$ws.Cells.Item(75, 2).Value = 0.198
$ws.Cells.Item(75, 3).Value = 0.312
$ws.Cells.Item(75, 4).Value = 0
$ws.Cells.Item(75, 5).Value = 75
$ws.Cells.Item(75, 6).Value = 74
$ws.Cells.Item(75, 7).Value = "01.10.202574"

$ws.Cells.Item(76, 1).Value = 45931.77083333334
$ws.Cells.Item(76, 2).Value = 0.15
$ws.Cells.Item(76, 3).Value = 0.27
$ws.Cells.Item(76, 4).Value = 0
$ws.Cells.Item(76, 5).Value = 75
$ws.Cells.Item(76, 6).Value = 75
$ws.Cells.Item(76, 7).Value = "01.10.202575"

$ws.Cells.Item(77, 1).Value = 45931.78125
$ws.Cells.Item(77, 2).Value = 0.224
$ws.Cells.Item(77, 3).Value = 0.1
$ws.Cells.Item(77, 4).Value = 0
$ws.Cells.Item(77, 5).Value = 75
$ws.Cells.Item(77, 6).Value = 76
$ws.Cells.Item(77, 7).Value = "01.10.202576"

$ws.Cells.Item(78, 1).Value = 45931.79166666666
$ws.Cells.Item(78, 2).Value = 0.255
$ws.Cells.Item(78, 3).Value = 0.427
$ws.Cells.Item(78, 4).Value = 0
$ws.Cells.Item(78, 5).Value = 75
$ws.Cells.Item(78, 6).Value = 77
$ws.Cells.Item(78, 7).Value = "01.10.202577"

$ws.Cells.Item(79, 1).Value = 45931.80208333334
$ws.Cells.Item(79, 2).Value = 0.109
$ws.Cells.Item(79, 3).Value = 0.125
$ws.Cells.Item(79, 4).Value = 0
$ws.Cells.Item(79, 5).Value = 75
$ws.Cells.Item(79, 6).Value = 78
$ws.Cells.Item(79, 7).Value = "01.10.202578"

$ws.Cells.Item(80, 1).Value = 45931.8125
$ws.Cells.Item(80, 2).Value = 0.033
$ws.Cells.Item(80, 3).Value = 7.878
$ws.Cells.Item(80, 4).Value = 0
$ws.Cells.Item(80, 5).Value = 75
$ws.Cells.Item(80, 6).Value = 79
$ws.Cells.Item(80, 7).Value = "01.10.202579"

$ws.Cells.Item(81, 1).Value = 45931.82291666666
$ws.Cells.Item(81, 2).Value = 0
$ws.Cells.Item(81, 3).Value = 5.211
$ws.Cells.Item(81, 4).Value = 0
$ws.Cells.Item(81, 5).Value = 75
$ws.Cells.Item(81, 6).Value = 80
$ws.Cells.Item(81, 7).Value = "01.10.202580"

$ws.Cells.Item(82, 1).Value = 45931.83333333334
$ws.Cells.Item(82, 2).Value = 1.609
$ws.Cells.Item(82, 3).Value = 2.557
$ws.Cells.Item(82, 4).Value = 0
$ws.Cells.Item(82, 5).Value = 75
$ws.Cells.Item(82, 6).Value = 81
$ws.Cells.Item(82, 7).Value = "01.10.202581"

$ws.Cells.Item(83, 1).Value = 45931.84375
$ws.Cells.Item(83, 2).Value = 0
$ws.Cells.Item(83, 3).Value = 5.151
$ws.Cells.Item(83, 4).Value = 0
$ws.Cells.Item(83, 5).Value = 75
$ws.Cells.Item(83, 6).Value = 82
$ws.Cells.Item(83, 7).Value = "01.10.202582"

$ws.Cells.Item(84, 1).Value = 45931.85416666666
$ws.Cells.Item(84, 2).Value = 4.137
$ws.Cells.Item(84, 3).Value = 0.212
$ws.Cells.Item(84, 4).Value = 0
$ws.Cells.Item(84, 5).Value = 65.75
$ws.Cells.Item(84, 6).Value = 83
$ws.Cells.Item(84, 7).Value = "01.10.202583"

$ws.Cells.Item(85, 1).Value = 45931.86458333334
$ws.Cells.Item(85, 2).Value = 14.197
$ws.Cells.Item(85, 3).Value = 0.026
$ws.Cells.Item(85, 4).Value = 0
$ws.Cells.Item(85, 5).Value = 65.75
$ws.Cells.Item(85, 6).Value = 84
$ws.Cells.Item(85, 7).Value = "01.10.202584"

$ws.Cells.Item(86, 1).Value = 45931.875
$ws.Cells.Item(86, 2).Value = 8.375
$ws.Cells.Item(86, 3).Value = 0.002
$ws.Cells.Item(86, 4).Value = 0
$ws.Cells.Item(86, 5).Value = 56.25
$ws.Cells.Item(86, 6).Value = 85
$ws.Cells.Item(86, 7).Value = "01.10.202585"

$ws.Cells.Item(87, 1).Value = 45931.88541666666
$ws.Cells.Item(87, 2).Value = 2.006
$ws.Cells.Item(87, 3).Value = 0.167
$ws.Cells.Item(87, 4).Value = 0
$ws.Cells.Item(87, 5).Value = 37.5
$ws.Cells.Item(87, 6).Value = 86
$ws.Cells.Item(87, 7).Value = "01.10.202586"

$ws.Cells.Item(88, 1).Value = 45931.89583333334
$ws.Cells.Item(88, 2).Value = 0.669
$ws.Cells.Item(88, 3).Value = 0.019
$ws.Cells.Item(88, 4).Value = 0
$ws.Cells.Item(88, 5).Value = 37.5
$ws.Cells.Item(88, 6).Value = 87
$ws.Cells.Item(88, 7).Value = "01.10.202587"

$ws.Cells.Item(89, 1).Value = 45931.90625
$ws.Cells.Item(89, 2).Value = 10.063
$ws.Cells.Item(89, 3).Value = 0.014
$ws.Cells.Item(89, 4).Value = 0
$ws.Cells.Item(89, 5).Value = 37.5
$ws.Cells.Item(89, 6).Value = 88
$ws.Cells.Item(89, 7).Value = "01.10.202588"

$ws.Cells.Item(90, 1).Value = 45931.91666666666
$ws.Cells.Item(90, 2).Value = 0
$ws.Cells.Item(90, 3).Value = 16.812
$ws.Cells.Item(90, 4).Value = 0
$ws.Cells.Item(90, 5).Value = 0
$ws.Cells.Item(90, 6).Value = 89
$ws.Cells.Item(90, 7).Value = "01.10.202589"

$ws.Cells.Item(91, 1).Value = 45931.92708333334
$ws.Cells.Item(91, 2).Value = 0
$ws.Cells.Item(91, 3).Value = 17.914
$ws.Cells.Item(91, 4).Value = 0
$ws.Cells.Item(91, 5).Value = 0
$ws.Cells.Item(91, 6).Value = 90
$ws.Cells.Item(91, 7).Value = "01.10.202590"

$ws.Cells.Item(92, 1).Value = 45931.9375
$ws.Cells.Item(92, 2).Value = 0.007
$ws.Cells.Item(92, 3).Value = 1.479
$ws.Cells.Item(92, 4).Value = 0
$ws.Cells.Item(92, 5).Value = 0
$ws.Cells.Item(92, 6).Value = 91
$ws.Cells.Item(92, 7).Value = "01.10.202591"

$ws.Cells.Item(93, 1).Value = 45931.94791666666
$ws.Cells.Item(93, 2).Value = 0.39
$ws.Cells.Item(93, 3).Value = 0.042
$ws.Cells.Item(93, 4).Value = 0
$ws.Cells.Item(93, 5).Value = 0
$ws.Cells.Item(93, 6).Value = 92
$ws.Cells.Item(93, 7).Value = "01.10.202592"

$ws.Cells.Item(94, 1).Value = 45931.95833333334
$ws.Cells.Item(94, 2).Value = 0.161
$ws.Cells.Item(94, 3).Value = 0.216
$ws.Cells.Item(94, 4).Value = 0
$ws.Cells.Item(94, 5).Value = 0
$ws.Cells.Item(94, 6).Value = 93
$ws.Cells.Item(94, 7).Value = "01.10.202593"

$ws.Cells.Item(95, 1).Value = 45931.96875
$ws.Cells.Item(95, 2).Value = 0.005
$ws.Cells.Item(95, 3).Value = 1.736
$ws.Cells.Item(95, 4).Value = 0
$ws.Cells.Item(95, 5).Value = 0
$ws.Cells.Item(95, 6).Value = 94
$ws.Cells.Item(95, 7).Value = "01.10.202594"

$ws.Cells.Item(96, 1).Value = 45931.97916666666
$ws.Cells.Item(96, 2).Value = 0.214
$ws.Cells.Item(96, 3).Value = 1.206
$ws.Cells.Item(96, 4).Value = 0
$ws.Cells.Item(96, 5).Value = 0
$ws.Cells.Item(96, 6).Value = 95
$ws.Cells.Item(96, 7).Value = "01.10.202595"

$ws.Cells.Item(97, 1).Value = 45931.98958333334
$ws.Cells.Item(97, 2).Value = 0
$ws.Cells.Item(97, 3).Value = 5.447
$ws.Cells.Item(97, 4).Value = 0
$ws.Cells.Item(97, 5).Value = 0
$ws.Cells.Item(97, 6).Value = 96
$ws.Cells.Item(97, 7).Value = "01.10.202596"

$ws.Cells.Item(98, 1).Value = 45932
$ws.Cells.Item(98, 2).Value = 0
$ws.Cells.Item(98, 3).Value = 23.253
$ws.Cells.Item(98, 4).Value = 0
$ws.Cells.Item(98, 5).Value = 0
$ws.Cells.Item(98, 6).Value = 1
$ws.Cells.Item(98, 7).Value = "02.10.20251"

$ws.Cells.Item(99, 1).Value = 45932
$ws.Cells.Item(99, 2).Value = 0
$ws.Cells.Item(99, 3).Value = 23.253
$ws.Cells.Item(99, 4).Value = 0
$ws.Cells.Item(99, 5).Value = 0
$ws.Cells.Item(99, 6).Value = 1
$ws.Cells.Item(99, 7).Value = "02.10.20251"

$ws.Cells.Item(100, 1).Value = 45932.01041666666
$ws.Cells.Item(100, 2).Value = 0
$ws.Cells.Item(100, 3).Value = 2.779
$ws.Cells.Item(100, 4).Value = 0
$ws.Cells.Item(100, 5).Value = 0
$ws.Cells.Item(100, 6).Value = 2
$ws.Cells.Item(100, 7).Value = "02.10.20252"

$ws.Cells.Item(101, 1).Value = 45932.01041666666
$ws.Cells.Item(101, 2).Value = 0
$ws.Cells.Item(101, 3).Value = 2.779
$ws.Cells.Item(101, 4).Value = 0
$ws.Cells.Item(101, 5).Value = 0
$ws.Cells.Item(101, 6).Value = 2
$ws.Cells.Item(101, 7).Value = "02.10.20252"

$ws.Cells.Item(102, 1).Value = 45932.02083333334
$ws.Cells.Item(102, 2).Value = 3.262
$ws.Cells.Item(102, 3).Value = 0
$ws.Cells.Item(102, 4).Value = 0
$ws.Cells.Item(102, 5).Value = 0
$ws.Cells.Item(102, 6).Value = 3
$ws.Cells.Item(102, 7).Value = "02.10.20253"

$ws.Cells.Item(103, 1).Value = 45932.02083333334
$ws.Cells.Item(103, 2).Value = 3.262
$ws.Cells.Item(103, 3).Value = 0
$ws.Cells.Item(103, 4).Value = 0
$ws.Cells.Item(103, 5).Value = 0
$ws.Cells.Item(103, 6).Value = 3
$ws.Cells.Item(103, 7).Value = "02.10.20253"

$ws.Cells.Item(104, 1).Value = 45932.03125
$ws.Cells.Item(104, 2).Value = 5.774
$ws.Cells.Item(104, 3).Value = 0
$ws.Cells.Item(104, 4).Value = 0
$ws.Cells.Item(104, 5).Value = 0
$ws.Cells.Item(104, 6).Value = 4
$ws.Cells.Item(104, 7).Value = "02.10.20254"

$ws.Cells.Item(105, 1).Value = 45932.03125
$ws.Cells.Item(105, 2).Value = 5.774
$ws.Cells.Item(105, 3).Value = 0
$ws.Cells.Item(105, 4).Value = 0
$ws.Cells.Item(105, 5).Value = 0
$ws.Cells.Item(105, 6).Value = 4
$ws.Cells.Item(105, 7).Value = "02.10.20254"

$ws.Cells.Item(106, 1).Value = 45932.04166666666
$ws.Cells.Item(106, 2).Value = 1.892
$ws.Cells.Item(106, 3).Value = 0
$ws.Cells.Item(106, 4).Value = 0
$ws.Cells.Item(106, 5).Value = 0
$ws.Cells.Item(106, 6).Value = 5
$ws.Cells.Item(106, 7).Value = "02.10.20255"

$ws.Cells.Item(107, 1).Value = 45932.04166666666
$ws.Cells.Item(107, 2).Value = 1.892
$ws.Cells.Item(107, 3).Value = 0
$ws.Cells.Item(107, 4).Value = 0
$ws.Cells.Item(107, 5).Value = 0
$ws.Cells.Item(107, 6).Value = 5
$ws.Cells.Item(107, 7).Value = "02.10.20255"

$ws.Cells.Item(108, 1).Value = 45932.05208333334
$ws.Cells.Item(108, 2).Value = 0.147
$ws.Cells.Item(108, 3).Value = 0.23
$ws.Cells.Item(108, 4).Value = 0
$ws.Cells.Item(108, 5).Value = 0
$ws.Cells.Item(108, 6).Value = 6
$ws.Cells.Item(108, 7).Value = "02.10.20256"

$ws.Cells.Item(109, 1).Value = 45932.05208333334
$ws.Cells.Item(109, 2).Value = 0.147
$ws.Cells.Item(109, 3).Value = 0.23
$ws.Cells.Item(109, 4).Value = 0
$ws.Cells.Item(109, 5).Value = 0
$ws.Cells.Item(109, 6).Value = 6
$ws.Cells.Item(109, 7).Value = "02.10.20256"

$ws.Cells.Item(110, 1).Value = 45932.0625
$ws.Cells.Item(110, 2).Value = 0
$ws.Cells.Item(110, 3).Value = 5.198
$ws.Cells.Item(110, 4).Value = 0
$ws.Cells.Item(110, 5).Value = 0
$ws.Cells.Item(110, 6).Value = 7
$ws.Cells.Item(110, 7).Value = "02.10.20257"

$ws.Cells.Item(111, 1).Value = 45932.0625
$ws.Cells.Item(111, 2).Value = 0
$ws.Cells.Item(111, 3).Value = 5.198
$ws.Cells.Item(111, 4).Value = 0
$ws.Cells.Item(111, 5).Value = 0
$ws.Cells.Item(111, 6).Value = 7
$ws.Cells.Item(111, 7).Value = "02.10.20257"

$ws.Cells.Item(112, 1).Value = 45932.07291666666
$ws.Cells.Item(112, 2).Value = 0
$ws.Cells.Item(112, 3).Value = 2.686
$ws.Cells.Item(112, 4).Value = 0
$ws.Cells.Item(112, 5).Value = 0
$ws.Cells.Item(112, 6).Value = 8
$ws.Cells.Item(112, 7).Value = "02.10.20258"

$ws.Cells.Item(113, 1).Value = 45932.07291666666
$ws.Cells.Item(113, 2).Value = 0
$ws.Cells.Item(113, 3).Value = 2.686
$ws.Cells.Item(113, 4).Value = 0
$ws.Cells.Item(113, 5).Value = 0
$ws.Cells.Item(113, 6).Value = 8
$ws.Cells.Item(113, 7).Value = "02.10.20258"

$ws.Cells.Item(114, 1).Value = 45932.08333333334
$ws.Cells.Item(114, 2).Value = 0
$ws.Cells.Item(114, 3).Value = 2.145
$ws.Cells.Item(114, 4).Value = 0
$ws.Cells.Item(114, 5).Value = 0
$ws.Cells.Item(114, 6).Value = 9
$ws.Cells.Item(114, 7).Value = "02.10.20259"

$ws.Cells.Item(115, 1).Value = 45932.09375
$ws.Cells.Item(115, 2).Value = 0.036
$ws.Cells.Item(115, 3).Value = 1.152
$ws.Cells.Item(115, 4).Value = 0
$ws.Cells.Item(115, 5).Value = 0
$ws.Cells.Item(115, 6).Value = 10
$ws.Cells.Item(115, 7).Value = "02.10.202510"

$ws.Cells.Item(116, 1).Value = 45932.10416666666
$ws.Cells.Item(116, 2).Value = 0
$ws.Cells.Item(116, 3).Value = 10.985
$ws.Cells.Item(116, 4).Value = 0
$ws.Cells.Item(116, 5).Value = 0
$ws.Cells.Item(116, 6).Value = 11
$ws.Cells.Item(116, 7).Value = "02.10.202511"

$ws.Cells.Item(117, 1).Value = 45932.11458333334
$ws.Cells.Item(117, 2).Value = 0.03
$ws.Cells.Item(117, 3).Value = 2.873
$ws.Cells.Item(117, 4).Value = 0
$ws.Cells.Item(117, 5).Value = 0
$ws.Cells.Item(117, 6).Value = 12
$ws.Cells.Item(117, 7).Value = "02.10.202512"

$ws.Cells.Item(118, 1).Value = 45932.125
$ws.Cells.Item(118, 2).Value = 1.723
$ws.Cells.Item(118, 3).Value = 0.43
$ws.Cells.Item(118, 4).Value = 0
$ws.Cells.Item(118, 5).Value = 0
$ws.Cells.Item(118, 6).Value = 13
$ws.Cells.Item(118, 7).Value = "02.10.202513"

$ws.Cells.Item(119, 1).Value = 45932.13541666666
$ws.Cells.Item(119, 2).Value = 12.738
$ws.Cells.Item(119, 3).Value = 0
$ws.Cells.Item(119, 4).Value = 0
$ws.Cells.Item(119, 5).Value = 0
$ws.Cells.Item(119, 6).Value = 14
$ws.Cells.Item(119, 7).Value = "02.10.202514"

$ws.Cells.Item(120, 1).Value = 45932.14583333334
$ws.Cells.Item(120, 2).Value = 1.541
$ws.Cells.Item(120, 3).Value = 0.125
$ws.Cells.Item(120, 4).Value = 0
$ws.Cells.Item(120, 5).Value = 0
$ws.Cells.Item(120, 6).Value = 15
$ws.Cells.Item(120, 7).Value = "02.10.202515"

$ws.Cells.Item(121, 1).Value = 45932.15625
$ws.Cells.Item(121, 2).Value = 22.776
$ws.Cells.Item(121, 3).Value = 0
$ws.Cells.Item(121, 4).Value = 0
$ws.Cells.Item(121, 5).Value = 0
$ws.Cells.Item(121, 6).Value = 16
$ws.Cells.Item(121, 7).Value = "02.10.202516"

$ws.Cells.Item(122, 1).Value = 45932.16666666666
$ws.Cells.Item(122, 2).Value = 10.996
$ws.Cells.Item(122, 3).Value = 0
$ws.Cells.Item(122, 4).Value = 0
$ws.Cells.Item(122, 5).Value = 0
$ws.Cells.Item(122, 6).Value = 17
$ws.Cells.Item(122, 7).Value = "02.10.202517"

$ws.Cells.Item(123, 1).Value = 45932.17708333334
$ws.Cells.Item(123, 2).Value = 8.347
$ws.Cells.Item(123, 3).Value = 0
$ws.Cells.Item(123, 4).Value = 0
$ws.Cells.Item(123, 5).Value = 0
$ws.Cells.Item(123, 6).Value = 18
$ws.Cells.Item(123, 7).Value = "02.10.202518"

$ws.Cells.Item(124, 1).Value = 45932.1875
$ws.Cells.Item(124, 2).Value = 9.102
$ws.Cells.Item(124, 3).Value = 0
$ws.Cells.Item(124, 4).Value = 0
$ws.Cells.Item(124, 5).Value = 0
$ws.Cells.Item(124, 6).Value = 19
$ws.Cells.Item(124, 7).Value = "02.10.202519"

$ws.Cells.Item(125, 1).Value = 45932.19791666666
$ws.Cells.Item(125, 2).Value = 4.434
$ws.Cells.Item(125, 3).Value = 0
$ws.Cells.Item(125, 4).Value = 0
$ws.Cells.Item(125, 5).Value = 0
$ws.Cells.Item(125, 6).Value = 20
$ws.Cells.Item(125, 7).Value = "02.10.202520"

$ws.Cells.Item(126, 1).Value = 45932.20833333334
$ws.Cells.Item(126, 2).Value = 12.237
$ws.Cells.Item(126, 3).Value = 0
$ws.Cells.Item(126, 4).Value = 0
$ws.Cells.Item(126, 5).Value = 0
$ws.Cells.Item(126, 6).Value = 21
$ws.Cells.Item(126, 7).Value = "02.10.202521"

$ws.Cells.Item(127, 1).Value = 45932.21875
$ws.Cells.Item(127, 2).Value = 1.511
$ws.Cells.Item(127, 3).Value = 0.5629999999999999
$ws.Cells.Item(127, 4).Value = 0
$ws.Cells.Item(127, 5).Value = 0
$ws.Cells.Item(127, 6).Value = 22
$ws.Cells.Item(127, 7).Value = "02.10.202522"

$ws.Cells.Item(128, 1).Value = 45932.22916666666
$ws.Cells.Item(128, 2).Value = 0
$ws.Cells.Item(128, 3).Value = 18.05
$ws.Cells.Item(128, 4).Value = 0
$ws.Cells.Item(128, 5).Value = 0
$ws.Cells.Item(128, 6).Value = 23
$ws.Cells.Item(128, 7).Value = "02.10.202523"

$ws.Cells.Item(129, 1).Value = 45932.23958333334
$ws.Cells.Item(129, 2).Value = 0.333
$ws.Cells.Item(129, 3).Value = 3.134
$ws.Cells.Item(129, 4).Value = 0
$ws.Cells.Item(129, 5).Value = 0
$ws.Cells.Item(129, 6).Value = 24
$ws.Cells.Item(129, 7).Value = "02.10.202524"

$ws.Cells.Item(130, 1).Value = 45932.25
$ws.Cells.Item(130, 2).Value = 27.551
$ws.Cells.Item(130, 3).Value = 0
$ws.Cells.Item(130, 4).Value = 0
$ws.Cells.Item(130, 5).Value = 0
$ws.Cells.Item(130, 6).Value = 25
$ws.Cells.Item(130, 7).Value = "02.10.202525"

$ws.Cells.Item(131, 1).Value = 45932.26041666666
$ws.Cells.Item(131, 2).Value = 41.72
$ws.Cells.Item(131, 3).Value = 0
$ws.Cells.Item(131, 4).Value = 0
$ws.Cells.Item(131, 5).Value = 0
$ws.Cells.Item(131, 6).Value = 26
$ws.Cells.Item(131, 7).Value = "02.10.202526"

$ws.Cells.Item(132, 1).Value = 45932.27083333334
$ws.Cells.Item(132, 2).Value = 27.496
$ws.Cells.Item(132, 3).Value = 0
$ws.Cells.Item(132, 4).Value = 0
$ws.Cells.Item(132, 5).Value = 0
$ws.Cells.Item(132, 6).Value = 27
$ws.Cells.Item(132, 7).Value = "02.10.202527"

$ws.Cells.Item(133, 1).Value = 45932.28125
$ws.Cells.Item(133, 2).Value = 27.566
$ws.Cells.Item(133, 3).Value = 0
$ws.Cells.Item(133, 4).Value = 0
$ws.Cells.Item(133, 5).Value = 0
$ws.Cells.Item(133, 6).Value = 28
$ws.Cells.Item(133, 7).Value = "02.10.202528"

$ws.Cells.Item(134, 1).Value = 45932.29166666666
$ws.Cells.Item(134, 2).Value = 50.987
$ws.Cells.Item(134, 3).Value = 0
$ws.Cells.Item(134, 4).Value = 0
$ws.Cells.Item(134, 5).Value = 0
$ws.Cells.Item(134, 6).Value = 29
$ws.Cells.Item(134, 7).Value = "02.10.202529"

$ws.Cells.Item(135, 1).Value = 45932.30208333334
$ws.Cells.Item(135, 2).Value = 47.114
$ws.Cells.Item(135, 3).Value = 0
$ws.Cells.Item(135, 4).Value = 14.5
$ws.Cells.Item(135, 5).Value = 0
$ws.Cells.Item(135, 6).Value = 30
$ws.Cells.Item(135, 7).Value = "02.10.202530"

$ws.Cells.Item(136, 1).Value = 45932.3125
$ws.Cells.Item(136, 2).Value = 37.56
$ws.Cells.Item(136, 3).Value = 0
$ws.Cells.Item(136, 4).Value = 32
$ws.Cells.Item(136, 5).Value = 0
$ws.Cells.Item(136, 6).Value = 31
$ws.Cells.Item(136, 7).Value = "02.10.202531"

$ws.Cells.Item(137, 1).Value = 45932.32291666666
$ws.Cells.Item(137, 2).Value = 18.966
$ws.Cells.Item(137, 3).Value = 0.021
$ws.Cells.Item(137, 4).Value = 50
$ws.Cells.Item(137, 5).Value = 0
$ws.Cells.Item(137, 6).Value = 32
$ws.Cells.Item(137, 7).Value = "02.10.202532"

$ws.Cells.Item(138, 1).Value = 45932.33333333334
$ws.Cells.Item(138, 2).Value = 20.81
$ws.Cells.Item(138, 3).Value = 0
$ws.Cells.Item(138, 4).Value = 50
$ws.Cells.Item(138, 5).Value = 0
$ws.Cells.Item(138, 6).Value = 33
$ws.Cells.Item(138, 7).Value = "02.10.202533"

$ws.Cells.Item(139, 1).Value = 45932.34375
$ws.Cells.Item(139, 2).Value = 10.298
$ws.Cells.Item(139, 3).Value = 0.001
$ws.Cells.Item(139, 4).Value = 41
$ws.Cells.Item(139, 5).Value = 0
$ws.Cells.Item(139, 6).Value = 34
$ws.Cells.Item(139, 7).Value = "02.10.202534"

$ws.Cells.Item(140, 1).Value = 45932.35416666666
$ws.Cells.Item(140, 2).Value = 13.979
$ws.Cells.Item(140, 3).Value = 0.001
$ws.Cells.Item(140, 4).Value = 41.75
$ws.Cells.Item(140, 5).Value = 0
$ws.Cells.Item(140, 6).Value = 35
$ws.Cells.Item(140, 7).Value = "02.10.202535"

$ws.Cells.Item(141, 1).Value = 45932.36458333334
$ws.Cells.Item(141, 2).Value = 0.522
$ws.Cells.Item(141, 3).Value = 0.268
$ws.Cells.Item(141, 4).Value = 41.75
$ws.Cells.Item(141, 5).Value = 0
$ws.Cells.Item(141, 6).Value = 36
$ws.Cells.Item(141, 7).Value = "02.10.202536"

$ws.Cells.Item(142, 1).Value = 45932.375
$ws.Cells.Item(142, 2).Value = 0.887
$ws.Cells.Item(142, 3).Value = 0.176
$ws.Cells.Item(142, 4).Value = 50
$ws.Cells.Item(142, 5).Value = 0
$ws.Cells.Item(142, 6).Value = 37
$ws.Cells.Item(142, 7).Value = "02.10.202537"

$ws.Cells.Item(143, 1).Value = 45932.38541666666
$ws.Cells.Item(143, 2).Value = 0
$ws.Cells.Item(143, 3).Value = 8.837999999999999
$ws.Cells.Item(143, 4).Value = 75
$ws.Cells.Item(143, 5).Value = 0
$ws.Cells.Item(143, 6).Value = 38
$ws.Cells.Item(143, 7).Value = "02.10.202538"

$ws.Cells.Item(144, 1).Value = 45932.39583333334
$ws.Cells.Item(144, 2).Value = 0.026
$ws.Cells.Item(144, 3).Value = 5.754
$ws.Cells.Item(144, 4).Value = 64.5
$ws.Cells.Item(144, 5).Value = 0
$ws.Cells.Item(144, 6).Value = 39
$ws.Cells.Item(144, 7).Value = "02.10.202539"

$ws.Cells.Item(145, 1).Value = 45932.40625
$ws.Cells.Item(145, 2).Value = 0.226
$ws.Cells.Item(145, 3).Value = 0.9419999999999999
$ws.Cells.Item(145, 4).Value = 75
$ws.Cells.Item(145, 5).Value = 0
$ws.Cells.Item(145, 6).Value = 40
$ws.Cells.Item(145, 7).Value = "02.10.202540"

$ws.Cells.Item(146, 1).Value = 45932.42708333334
$ws.Cells.Item(146, 2).Value = 0
$ws.Cells.Item(146, 3).Value = 0
$ws.Cells.Item(146, 4).Value = 87.5
$ws.Cells.Item(146, 5).Value = 0
$ws.Cells.Item(146, 6).Value = 42
$ws.Cells.Item(146, 7).Value = "02.10.202542"

for ($r = 134; $r -le 146; $r++) {
    $ws.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
}
